$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'leg chaffing protection',
    'leg chaffing stick',
    'leg compresion',
    'leg compression',
    'leg compression basketball',
    'leg compression device',
    'leg compression devices',
    'leg compression for women',
    'leg compression full',
    'leg compression garment',
    'leg compression leggings',
    'leg compression machines',
    'leg compression pants',
    'leg compression pants men',
    'leg compression recovery',
    'leg compression running',
    'leg compression therapy',
    'leg compression tights',
    'leg conditioning',
    'leg drop panel',
    'leg muscle',
    'leg muscles',
    'leg pack',
    'leg pads basketball',
    'leg squeeze machine',
    'leg straightening brace',
    'leg stretch',
    'leg support for running',
    'leg support leggings',
    'leg supports',
    'leg tights for women',
    'leg warmer bike',
    'leg warmers ballet women',
    'leg warmers compression',
    'leg warmers cycling women',
    'leg warmers for basketball',
    'leg workout compression pants women',
    'leggens for women',
    'leggin shorts',
    'leggin victoria secret for women',
    'legging bands',
    'legging capri',
    'legging capris for women',
    'legging compression women',
    'legging design',
    'legging for basketball boys',
    'legging for men workout',
    'legging for running',
    'legging for skiing',
    'legging gym',
    'legging high waist',
    'legging knee length',
    'legging multipack',
    'legging nike women',
    'legging package',
    'legging pads',
    'legging patches for holes',
    'legging running',
    'legging running women',
    'legging short',
    'legging shorts for women high waist',
    'legging shorts women',
    'legging skirt',
    'legging skirt women',
    'legging tights',
    'legging underwear women',
    'legging victoria secret',
    'legging winter men',
    'legging with back pockets',
    'legging with designs',
    'legging with holes',
    'legging with patch knees',
    'legging with skirt',
    'legging women nike',
    'legginga for women',
    'leggings 3 pack',
    'leggings 360',
    'leggings 7 pack',
    'leggings 7 pack women',
    'leggings active wear for women',
    'leggings and tights',
    'leggings athleta',
    'leggings athletic',
    'leggings basketball',
    'leggings black capri',
    'leggings capri',
    'leggings capri for women',
    'leggings capri length',
    'leggings capris',
    'leggings cold weather womens',
    'leggings compression',
    'leggings compression for women',
    'leggings compression women',
    'leggings cool design',
    'leggings cortos',
    'leggings cortos for women',
    'leggings cycling',
    'leggings exercise',
    'leggings fitted',
    'leggings for basketball'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
